$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44511
$ws.Range("K2").Value = 'Start Ruby'
$ws.Range("L2").Value = 'Primera'
$ws.Range("M2").Value = 24
$ws.Range("N2").Value = 140000
$ws.Range("O2").Value = 150000
$ws.Range("P2").Value = 145000
$ws.Range("Q2").Value = '$/bins (350 kilos)'
$ws.Range("R2").Value = 'Región Metropolitana'
$ws.Range("S2").Value = 414
$ws.Range("T2").Value = 350

# Row 3
$ws.Range("D3").Value = 44196
$ws.Range("K3").Value = 'Red Blush'
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 12
$ws.Range("N3").Value = 130000
$ws.Range("O3").Value = 130000
$ws.Range("P3").Value = 130000
$ws.Range("Q3").Value = '$/bins (350 kilos)'
$ws.Range("R3").Value = 'Provincia de Limarí'
$ws.Range("S3").Value = 371
$ws.Range("T3").Value = 350

# Row 4
$ws.Range("D4").Value = 44193
$ws.Range("K4").Value = 'Start Ruby'
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 8
$ws.Range("N4").Value = 150000
$ws.Range("O4").Value = 150000
$ws.Range("P4").Value = 150000
$ws.Range("Q4").Value = '$/bins (350 kilos)'
$ws.Range("R4").Value = 'Región Metropolitana'
$ws.Range("S4").Value = 429
$ws.Range("T4").Value = 350

# Row 5
$ws.Range("D5").Value = 44298
$ws.Range("K5").Value = 'Start Ruby'
$ws.Range("L5").Value = 'Especial'
$ws.Range("M5").Value = 15
$ws.Range("N5").Value = 450000
$ws.Range("O5").Value = 450000
$ws.Range("P5").Value = 450000
$ws.Range("Q5").Value = '$/bins (350 kilos)'
$ws.Range("R5").Value = 'Región Metropolitana'
$ws.Range("S5").Value = 1286
$ws.Range("T5").Value = 350

# Row 6
$ws.Range("D6").Value = 44298
$ws.Range("K6").Value = 'Start Ruby'
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 20
$ws.Range("N6").Value = 430000
$ws.Range("O6").Value = 430000
$ws.Range("P6").Value = 430000
$ws.Range("Q6").Value = '$/bins (350 kilos)'
$ws.Range("R6").Value = 'Región Metropolitana'
$ws.Range("S6").Value = 1229
$ws.Range("T6").Value = 350

# Row 7
$ws.Range("D7").Value = 44167
$ws.Range("K7").Value = 'Start Ruby'
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 140
$ws.Range("N7").Value = 9800
$ws.Range("O7").Value = 9800
$ws.Range("P7").Value = 9800
$ws.Range("Q7").Value = '$/caja 14 kilos empedrada'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value = 700
$ws.Range("T7").Value = 14

# Row 8
$ws.Range("D8").Value = 44308
$ws.Range("K8").Value = 'Start Ruby'
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 20
$ws.Range("N8").Value = 280000
$ws.Range("O8").Value = 280000
$ws.Range("P8").Value = 280000
$ws.Range("Q8").Value = '$/bins (350 kilos)'
$ws.Range("R8").Value = 'Región Metropolitana'
$ws.Range("S8").Value = 800
$ws.Range("T8").Value = 350

# Row 9
$ws.Range("D9").Value = 44208
$ws.Range("K9").Value = 'Start Ruby'
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 16
$ws.Range("N9").Value = 180000
$ws.Range("O9").Value = 180000
$ws.Range("P9").Value = 180000
$ws.Range("Q9").Value = '$/bins (350 kilos)'
$ws.Range("R9").Value = 'Región Metropolitana'
$ws.Range("S9").Value = 514
$ws.Range("T9").Value = 350

# Row 10
$ws.Range("D10").Value = 44505
$ws.Range("K10").Value = 'Start Ruby'
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 15
$ws.Range("N10").Value = 150000
$ws.Range("O10").Value = 150000
$ws.Range("P10").Value = 150000
$ws.Range("Q10").Value = '$/bins (350 kilos)'
$ws.Range("R10").Value = 'Provincia de Quillota'
$ws.Range("S10").Value = 429
$ws.Range("T10").Value = 350

# Row 11
$ws.Range("D11").Value = 44189
$ws.Range("K11").Value = 'Start Ruby'
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 16
$ws.Range("N11").Value = 150000
$ws.Range("O11").Value = 150000
$ws.Range("P11").Value = 150000
$ws.Range("Q11").Value = '$/bins (350 kilos)'
$ws.Range("R11").Value = 'Provincia de Limarí'
$ws.Range("S11").Value = 429
$ws.Range("T11").Value = 350

# Row 12
$ws.Range("D12").Value = 44356
$ws.Range("K12").Value = 'Start Ruby'
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 24
$ws.Range("N12").Value = 200000
$ws.Range("O12").Value = 230000
$ws.Range("P12").Value = 215000
$ws.Range("Q12").Value = '$/bins (350 kilos)'
$ws.Range("R12").Value = 'Región Metropolitana'
$ws.Range("S12").Value = 614
$ws.Range("T12").Value = 350

# Row 13
$ws.Range("D13").Value = 44446
$ws.Range("K13").Value = 'Start Ruby'
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 14
$ws.Range("N13").Value = 150000
$ws.Range("O13").Value = 160000
$ws.Range("P13").Value = 155000
$ws.Range("Q13").Value = '$/bins (350 kilos)'
$ws.Range("R13").Value = 'Región de O''Higgins'
$ws.Range("S13").Value = 443
$ws.Range("T13").Value = 350

# Row 14
$ws.Range("D14").Value = 44363
$ws.Range("K14").Value = 'Start Ruby'
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 20
$ws.Range("N14").Value = 200000
$ws.Range("O14").Value = 230000
$ws.Range("P14").Value = 215000
$ws.Range("Q14").Value = '$/bins (350 kilos)'
$ws.Range("R14").Value = 'Provincia de Limarí'
$ws.Range("S14").Value = 614
$ws.Range("T14").Value = 350

# Row 15
$ws.Range("D15").Value = 44195
$ws.Range("K15").Value = 'Start Ruby'
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 20
$ws.Range("N15").Value = 200000
$ws.Range("O15").Value = 210000
$ws.Range("P15").Value = 206000
$ws.Range("Q15").Value = '$/bins (350 kilos)'
$ws.Range("R15").Value = 'Región de O''Higgins'
$ws.Range("S15").Value = 589
$ws.Range("T15").Value = 350

# Row 16
$ws.Range("D16").Value = 44376
$ws.Range("K16").Value = 'Start Ruby'
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 20
$ws.Range("N16").Value = 180000
$ws.Range("O16").Value = 180000
$ws.Range("P16").Value = 180000
$ws.Range("Q16").Value = '$/bins (350 kilos)'
$ws.Range("R16").Value = 'Hijuelas'
$ws.Range("S16").Value = 514
$ws.Range("T16").Value = 350

# Row 17
$ws.Range("D17").Value = 44376
$ws.Range("K17").Value = 'Start Ruby'
$ws.Range("L17").Value = 'Segunda'
$ws.Range("M17").Value = 16
$ws.Range("N17").Value = 150000
$ws.Range("O17").Value = 150000
$ws.Range("P17").Value = 150000
$ws.Range("Q17").Value = '$/bins (350 kilos)'
$ws.Range("R17").Value = 'Provincia de Limarí'
$ws.Range("S17").Value = 429
$ws.Range("T17").Value = 350

# Row 18
$ws.Range("D18").Value = 44312
$ws.Range("K18").Value = 'Start Ruby'
$ws.Range("L18").Value = 'Segunda'
$ws.Range("M18").Value = 10
$ws.Range("N18").Value = 330000
$ws.Range("O18").Value = 330000
$ws.Range("P18").Value = 330000
$ws.Range("Q18").Value = '$/bins (350 kilos)'
$ws.Range("R18").Value = 'Región Metropolitana'
$ws.Range("S18").Value = 943
$ws.Range("T18").Value = 350

# Row 19
$ws.Range("D19").Value = 44389
$ws.Range("K19").Value = 'Start Ruby'
$ws.Range("L19").Value = 'Especial'
$ws.Range("M19").Value = 18
$ws.Range("N19").Value = 200000
$ws.Range("O19").Value = 200000
$ws.Range("P19").Value = 200000
$ws.Range("Q19").Value = '$/bins (350 kilos)'
$ws.Range("R19").Value = 'Provincia de Quillota'
$ws.Range("S19").Value = 571
$ws.Range("T19").Value = 350

# Row 20
$ws.Range("D20").Value = 44400
$ws.Range("K20").Value = 'Start Ruby'
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 140
$ws.Range("N20").Value = 9800
$ws.Range("O20").Value = 9800
$ws.Range("P20").Value = 9800
$ws.Range("Q20").Value = '$/caja 14 kilos empedrada'
$ws.Range("R20").Value = 'Región de O''Higgins'
$ws.Range("S20").Value = 700
$ws.Range("T20").Value = 14

# Row 21
$ws.Range("D21").Value = 44309
$ws.Range("K21").Value = 'Start Ruby'
$ws.Range("L21").Value = 'Primera'
$ws.Range("M21").Value = 16
$ws.Range("N21").Value = 350000
$ws.Range("O21").Value = 350000
$ws.Range("P21").Value = 350000
$ws.Range("Q21").Value = '$/bins (350 kilos)'
$ws.Range("R21").Value = 'Región Metropolitana'
$ws.Range("S21").Value = 1000
$ws.Range("T21").Value = 350

# Row 22
$ws.Range("D22").Value = 44201
$ws.Range("K22").Value = 'Start Ruby'
$ws.Range("L22").Value = 'Especial'
$ws.Range("M22").Value = 8
$ws.Range("N22").Value = 200000
$ws.Range("O22").Value = 200000
$ws.Range("P22").Value = 200000
$ws.Range("Q22").Value = '$/bins (350 kilos)'
$ws.Range("R22").Value = 'Región de O''Higgins'
$ws.Range("S22").Value = 571
$ws.Range("T22").Value = 350

# Row 23
$ws.Range("D23").Value = 44201
$ws.Range("K23").Value = 'Start Ruby'
$ws.Range("L23").Value = 'Primera'
$ws.Range("M23").Value = 16
$ws.Range("N23").Value = 170000
$ws.Range("O23").Value = 170000
$ws.Range("P23").Value = 170000
$ws.Range("Q23").Value = '$/bins (350 kilos)'
$ws.Range("R23").Value = 'Región de O''Higgins'
$ws.Range("S23").Value = 486
$ws.Range("T23").Value = 350
